$p = $ppt.ActivePresentation
$s = $p.Slides
$s | Get-Member | Out-String | Write-Host
